$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("FEINmismatch")
$ws1.Range("B2").Value = "Sat Feb 17 22:19:14 EST 2024"
$ws1.Range("B3").Value = "Sat Feb 17 22:19:29 EST 2024"
$ws1.Range("B4").Value = "Sat Feb 17 22:19:40 EST 2024"
$ws1.Range("B5").Value = "Sat Feb 17 22:19:52 EST 2024"
$ws1.Range("B6").Value = "Sat Feb 17 22:20:03 EST 2024"
$ws1.Range("B7").Value = "Sat Feb 17 22:20:14 EST 2024"
$ws1.Range("B8").Value = "Sat Feb 17 22:20:25 EST 2024"
$ws1.Range("B9").Value = "Sat Feb 17 22:20:36 EST 2024"
$ws1.Range("B10").Value = "Sat Feb 17 22:20:47 EST 2024"
$ws1.Range("B13").Value = "Sat Feb 17 22:20:59 EST 2024"
$ws1.Range("B14").Value = "Sat Feb 17 22:21:10 EST 2024"
$ws1.Range("B15").Value = "Sat Feb 17 22:21:21 EST 2024"
$ws1.Range("B16").Value = "Sat Feb 17 22:21:32 EST 2024"
$ws1.Range("B17").Value = "Sat Feb 17 22:21:42 EST 2024"
$ws1.Range("B18").Value = "Sat Feb 17 22:21:53 EST 2024"
$ws1.Range("B19").Value = "Sat Feb 17 22:22:04 EST 2024"
$ws1.Range("B20").Value = "Sat Feb 17 22:22:15 EST 2024"
$ws1.Range("B21").Value = "Sat Feb 17 22:22:26 EST 2024"
$ws1.Range("B22").Value = "Sat Feb 17 22:22:37 EST 2024"
$ws1.Range("B23").Value = "Sat Feb 17 22:22:48 EST 2024"
$ws1.Range("B24").Value = "Sat Feb 17 22:22:58 EST 2024"
$ws1.Range("B25").Value = "Sat Feb 17 22:23:09 EST 2024"
$ws1.Range("B26").Value = "Sat Feb 17 22:23:20 EST 2024"
$ws1.Range("B27").Value = "Sat Feb 17 22:23:31 EST 2024"
$ws1.Range("B28").Value = "Sat Feb 17 22:23:42 EST 2024"
$ws1.Range("B29").Value = "Sat Feb 17 22:23:53 EST 2024"
$ws1.Range("B30").Value = "Sat Feb 17 22:24:04 EST 2024"

$ws2 = $wb.Worksheets.Item("FEINSSNmismatch")
$ws2.Range("B2").Value = "Sat Feb 17 22:24:15 EST 2024"
$ws2.Range("B3").Value = "Sat Feb 17 22:24:26 EST 2024"
$ws2.Range("B4").Value = "Sat Feb 17 22:24:36 EST 2024"
$ws2.Range("B5").Value = "Sat Feb 17 22:24:47 EST 2024"
$ws2.Range("B6").Value = "Sat Feb 17 22:24:57 EST 2024"
$ws2.Range("B7").Value = "Sat Feb 17 22:25:08 EST 2024"
$ws2.Range("B8").Value = "Sat Feb 17 22:25:18 EST 2024"
$ws2.Range("B9").Value = "Sat Feb 17 22:25:29 EST 2024"
$ws2.Range("B14").Value = "Sat Feb 17 22:25:40 EST 2024"
$ws2.Range("B15").Value = "Sat Feb 17 22:25:50 EST 2024"
$ws2.Range("B16").Value = "Sat Feb 17 22:26:01 EST 2024"
$ws2.Range("B17").Value = "Sat Feb 17 22:26:12 EST 2024"
$ws2.Range("B18").Value = "Sat Feb 17 22:26:22 EST 2024"
$ws2.Range("B19").Value = "Sat Feb 17 22:26:33 EST 2024"
